$wb = $excel.ActiveWorkbook

# --- Sheet "y1" (sheet1) ---
$ws1 = $wb.Worksheets.Item("y1")
$ws1.Range("C2").Value = 0
$ws1.Range("D4").Value = 0
$ws1.Range("F4").Value = 1
$ws1.Range("A5").Value = 0
$ws1.Range("E8").Value = 1
$ws1.Range("F9").Value = 0

# --- Sheet "y2" (sheet2) ---
$ws2 = $wb.Worksheets.Item("y2")
$ws2.Range("B5").Value = 1
$ws2.Range("C7").Value = 1
$ws2.Range("A9").Value = 1

# --- Sheet "y3" (sheet3) ---
$ws3 = $wb.Worksheets.Item("y3")
$ws3.Range("E2").Value = 1
$ws3.Range("B3").Value = 0
$ws3.Range("F3").Value = 1
$ws3.Range("A6").Value = 0
$ws3.Range("D6").Value = 1
$ws3.Range("C7").Value = 0
$ws3.Range("E8").Value = 0
